$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of timetable data (Course Code, TYPE, Day, Start Time, End Time)
$ws.Range("A13").Value = "SC2207"
$ws.Range("B13").Value = "LEC/STUDIO"
$ws.Range("C13").Value = "THU"
$ws.Range("D13").Value = "16:30"
$ws.Range("E13").Value = "17:20"

$ws.Range("A14").Value = "SC2207"
$ws.Range("B14").Value = "LEC/STUDIO"
$ws.Range("C14").Value = "TUE"
$ws.Range("D14").Value = "16:30"
$ws.Range("E14").Value = "17:20"

$ws.Range("C15").Select()
